$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = 44315
$ws.Cells.Item(2, 13).Value = 45

# Row 5
$ws.Cells.Item(5, 4).Value = 44314
$ws.Cells.Item(5, 13).Value = 47
$ws.Cells.Item(5, 14).Value = 9000
$ws.Cells.Item(5, 15).Value = 9000
$ws.Cells.Item(5, 16).Value = 9000
$ws.Cells.Item(5, 19).Value = 900

# Row 6
$ws.Cells.Item(6, 4).Value = 44329
$ws.Cells.Item(6, 12).Value = 'Primera'
$ws.Cells.Item(6, 13).Value = 56
$ws.Cells.Item(6, 14).Value = 9000
$ws.Cells.Item(6, 15).Value = 9000
$ws.Cells.Item(6, 16).Value = 9000
$ws.Cells.Item(6, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(6, 19).Value = 900

# Row 7
$ws.Cells.Item(7, 4).Value = 44329
$ws.Cells.Item(7, 12).Value = 'Segunda'
$ws.Cells.Item(7, 13).Value = 50
$ws.Cells.Item(7, 14).Value = 8000
$ws.Cells.Item(7, 15).Value = 8000
$ws.Cells.Item(7, 16).Value = 8000
$ws.Cells.Item(7, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(7, 19).Value = 800

# Row 8
$ws.Cells.Item(8, 4).Value = 44302
$ws.Cells.Item(8, 12).Value = 'Primera'
$ws.Cells.Item(8, 13).Value = 45
$ws.Cells.Item(8, 14).Value = 10000
$ws.Cells.Item(8, 15).Value = 10000
$ws.Cells.Item(8, 16).Value = 10000
$ws.Cells.Item(8, 19).Value = 1000

# Row 9
$ws.Cells.Item(9, 4).Value = 44319
$ws.Cells.Item(9, 13).Value = 68

# Row 10
$ws.Cells.Item(10, 4).Value = 44319
$ws.Cells.Item(10, 13).Value = 57
$ws.Cells.Item(10, 14).Value = 8000
$ws.Cells.Item(10, 15).Value = 8000
$ws.Cells.Item(10, 16).Value = 8000
$ws.Cells.Item(10, 19).Value = 800

# Row 11
$ws.Cells.Item(11, 4).Value = 44322
$ws.Cells.Item(11, 13).Value = 56

# Row 12
$ws.Cells.Item(12, 4).Value = 44322
$ws.Cells.Item(12, 13).Value = 40

# Row 13
$ws.Cells.Item(13, 4).Value = 44333
$ws.Cells.Item(13, 12).Value = 'Especial'
$ws.Cells.Item(13, 13).Value = 58

# Row 14
$ws.Cells.Item(14, 4).Value = 44333
$ws.Cells.Item(14, 13).Value = 65
$ws.Cells.Item(14, 14).Value = 9000
$ws.Cells.Item(14, 15).Value = 9000
$ws.Cells.Item(14, 16).Value = 9000
$ws.Cells.Item(14, 19).Value = 900

# Row 15
$ws.Cells.Item(15, 4).Value = 44333
$ws.Cells.Item(15, 13).Value = 60
$ws.Cells.Item(15, 14).Value = 8000
$ws.Cells.Item(15, 15).Value = 8000
$ws.Cells.Item(15, 16).Value = 8000
$ws.Cells.Item(15, 19).Value = 800

# Row 16
$ws.Cells.Item(16, 4).Value = 44321
$ws.Cells.Item(16, 12).Value = 'Primera'
$ws.Cells.Item(16, 13).Value = 58
$ws.Cells.Item(16, 14).Value = 9000
$ws.Cells.Item(16, 15).Value = 9000
$ws.Cells.Item(16, 16).Value = 9000
$ws.Cells.Item(16, 18).Value = 'Provincia de Quillota'
$ws.Cells.Item(16, 19).Value = 900

# Row 17
$ws.Cells.Item(17, 4).Value = 44307
$ws.Cells.Item(17, 13).Value = 40
$ws.Cells.Item(17, 14).Value = 10000
$ws.Cells.Item(17, 15).Value = 10000
$ws.Cells.Item(17, 16).Value = 10000
$ws.Cells.Item(17, 18).Value = 'Provincia de Quillota'
$ws.Cells.Item(17, 19).Value = 1000

# Row 18
$ws.Cells.Item(18, 4).Value = 44328
$ws.Cells.Item(18, 12).Value = 'Primera'
$ws.Cells.Item(18, 13).Value = 45
$ws.Cells.Item(18, 18).Value = 'Provincia de Quillota'

# Row 19
$ws.Cells.Item(19, 4).Value = 44328
$ws.Cells.Item(19, 12).Value = 'Segunda'
$ws.Cells.Item(19, 13).Value = 48
$ws.Cells.Item(19, 14).Value = 7000
$ws.Cells.Item(19, 15).Value = 7000
$ws.Cells.Item(19, 16).Value = 7000
$ws.Cells.Item(19, 19).Value = 700

# Row 20
$ws.Cells.Item(20, 4).Value = 44309
$ws.Cells.Item(20, 13).Value = 45

# Row 21
$ws.Cells.Item(21, 4).Value = 44343
$ws.Cells.Item(21, 12).Value = 'Especial'
$ws.Cells.Item(21, 13).Value = 47
$ws.Cells.Item(21, 14).Value = 10000
$ws.Cells.Item(21, 15).Value = 10000
$ws.Cells.Item(21, 16).Value = 10000
$ws.Cells.Item(21, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(21, 19).Value = 1000

# Row 22
$ws.Cells.Item(22, 4).Value = 44343
$ws.Cells.Item(22, 13).Value = 50
$ws.Cells.Item(22, 18).Value = 'Región Metropolitana'

# Row 23
$ws.Cells.Item(23, 4).Value = 44343
$ws.Cells.Item(23, 12).Value = 'Segunda'
$ws.Cells.Item(23, 13).Value = 58
$ws.Cells.Item(23, 14).Value = 8000
$ws.Cells.Item(23, 15).Value = 8000
$ws.Cells.Item(23, 16).Value = 8000
$ws.Cells.Item(23, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(23, 19).Value = 800

# Row 24
$ws.Cells.Item(24, 4).Value = 44312
$ws.Cells.Item(24, 13).Value = 48

# Row 25
$ws.Cells.Item(25, 4).Value = 44308

# Row 26
$ws.Cells.Item(26, 4).Value = 44308
$ws.Cells.Item(26, 12).Value = 'Segunda'
$ws.Cells.Item(26, 13).Value = 48
$ws.Cells.Item(26, 14).Value = 8000
$ws.Cells.Item(26, 15).Value = 8000
$ws.Cells.Item(26, 16).Value = 8000
$ws.Cells.Item(26, 18).Value = 'Provincia de Quillota'
$ws.Cells.Item(26, 19).Value = 800

# Row 27
$ws.Cells.Item(27, 4).Value = 44301
$ws.Cells.Item(27, 12).Value = 'Primera'
$ws.Cells.Item(27, 13).Value = 45
$ws.Cells.Item(27, 14).Value = 10000
$ws.Cells.Item(27, 15).Value = 10000
$ws.Cells.Item(27, 16).Value = 10000
$ws.Cells.Item(27, 18).Value = 'Provincia de Quillota'
$ws.Cells.Item(27, 19).Value = 1000

# Row 28
$ws.Cells.Item(28, 4).Value = 44306
$ws.Cells.Item(28, 13).Value = 45

# Row 29
$ws.Cells.Item(29, 4).Value = 44323
$ws.Cells.Item(29, 12).Value = 'Primera'
$ws.Cells.Item(29, 13).Value = 60
$ws.Cells.Item(29, 14).Value = 10000
$ws.Cells.Item(29, 15).Value = 10000
$ws.Cells.Item(29, 16).Value = 10000
$ws.Cells.Item(29, 19).Value = 1000

# Row 30
$ws.Cells.Item(30, 4).Value = 44323
$ws.Cells.Item(30, 12).Value = 'Segunda'
$ws.Cells.Item(30, 13).Value = 50
